{"js": "// The target paragraph currently reads (across two runs):\n//   run1: \"\u00e0 realiza\u00e7\u00e3o de vendas\"\n//   run2: \" ser feita exclusivamente na loja f\u00edsica\"\n// After the edit it must read (again across two runs):\n//   run1: \"\u00e0 \"\n//   run2: \"falta de investimentos em novas tecnologias para melhorar o\n//          atendimento captando novos clientes e fidelizando os atuais  \"\n//\n// We perform two separate search/replace operations, each confined to a\n// single existing run, so Word keeps the run split (and each run's own\n// formatting/rsid) intact instead of merging everything into one run.\n\nconst body = context.document.body;\n\n// 1) Shrink the first run down to \"\u00e0 \" (trailing space preserved).\nconst firstRunResults = body.search(\"\u00e0 realiza\u00e7\u00e3o de vendas\", { matchCase: true, matchWholeWord: false });\nfirstRunResults.load(\"items\");\nawait context.sync();\n\nif (firstRunResults.items.length > 0) {\n  firstRunResults.items[0].insertText(\"\u00e0 \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Replace the text of the second run with the new sentence.\nconst secondRunResults = body.search(\" ser feita exclusivamente na loja f\u00edsica\", { matchCase: true, matchWholeWord: false });\nsecondRunResults.load(\"items\");\nawait context.sync();\n\nif (secondRunResults.items.length > 0) {\n  secondRunResults.items[0].insertText(\n    \"falta de investimentos em novas tecnologias para melhorar o atendimento captando novos clientes e fidelizando os atuais  \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# The \"affeta...devido\" table row currently reads (across two runs):\n#   run1: \"\u00e0 realiza\u00e7\u00e3o de vendas\"\n#   run2: \" ser feita exclusivamente na loja f\u00edsica\"\n# After the edit it must read (again across two runs):\n#   run1: \"\u00e0 \"\n#   run2: \"falta de investimentos em novas tecnologias para melhorar o\n#          atendimento captando novos clientes e fidelizando os atuais  \"\n#\n# Word's Range.Text setter normally re-merges adjacent runs that end up with\n# identical formatting, so a naive pair of replacements would collapse both\n# runs into a single <w:r>. To keep the document's run structure faithful to\n# the target (two separate runs), run2 is replaced while briefly toggled to\n# Bold so Word is forced to treat it as a distinct run; once its text is in\n# place we flip Bold back off, leaving the formatting unchanged from the\n# original (Arial, no bold) but the run split intact.\n\n$d = $word.ActiveDocument\n\n$oldRun1Text = \"\u00e0 realiza\u00e7\u00e3o de vendas\"\n$newRun1Text = \"\u00e0 \"\n\n$oldRun2Text = \" ser feita exclusivamente na loja f\u00edsica\"\n$newRun2Text = \"falta de investimentos em novas tecnologias para melhorar o atendimento captando novos clientes e fidelizando os atuais  \"\n\n# --- Run 1: plain in-place text replacement, formatting untouched. ---\n$run1 = $d.Content\n$run1.Find.ClearFormatting()\n$run1.Find.Text = $oldRun1Text\n$run1.Find.Execute() | Out-Null\n$run1.Text = $newRun1Text\n\n# --- Run 2: replace text, forcing (then clearing) Bold to preserve the run split. ---\n$run2 = $d.Content\n$run2.Find.ClearFormatting()\n$run2.Find.Text = $oldRun2Text\n$run2.Find.Execute() | Out-Null\n$run2Start = $run2.Start\n$run2.Font.Bold = 1\n$run2.Text = $newRun2Text\n\n$run2Fixed = $d.Range($run2Start, $run2Start + $newRun2Text.Length)\n$run2Fixed.Font.Bold = 0\n"}
